$d = $word.ActiveDocument

# Locate the last paragraph in the document body (the one ending in
# "...jotta tiedot olisivat helpompia lukea.") and position a collapsed
# range right after it so we can append the new content.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
[void]$r.InsertParagraphAfter()

# InsertParagraphAfter created a new, empty paragraph at the end of the
# document (before the sectPr). Grab its range and inject the new
# "Modulaarisuus" heading plus the following body paragraph as raw OOXML
# so we get the exact run/proofErr structure the edit calls for.
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$nr = $newPara.Range

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
  </w:pPr>
  <w:r>
    <w:t>Modulaarisuus</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">Paransin ohjelman dynaamisuutta. Tein navigaatio </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>generaation</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> täysin dynaamiseksi käyttämällä kansio rakennetta, josta se hakee kansioiden alla olevat </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>usercontrol</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> tiedostot. Näiden perusteella ohjelma luo navigaatio palkin. Myös muutin asetuksia hakemaan normi sivua varten olevat sivut dynaamisesti tämän muutoksen perusteella. </w:t>
  </w:r>
</w:p>
"@

[void]$nr.InsertXML($xml)
